$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H62").Value = 2456.5
$ws.Range("I62").Value = 3296.3333
$ws.Range("K62").Value = 3296.3333
$ws.Range("M62").Value = -2672.3333
$ws.Range("H65").Value = 2456.5
$ws.Range("I65").Value = 3296.3333
$ws.Range("K65").Value = 16481.6665
$ws.Range("M65").Value = -13361.6665
$ws.Range("H106").Value = 2048.5715
$ws.Range("I106").Value = 2016.5264
$ws.Range("K106").Value = 2016.5264
$ws.Range("M106").Value = -1385.5264
$ws.Range("H132").Value = 7624.722
$ws.Range("I132").Value = 4187.0835
$ws.Range("K132").Value = 12561.2505
$ws.Range("M132").Value = -10031.2505
$ws.Range("H138").Value = 1983.5253
$ws.Range("I138").Value = 1405.0667
$ws.Range("J138").Value = 2086.8215
$ws.Range("K138").Value = 4215.2001
$ws.Range("L138").Value = 6260.4645
$ws.Range("M138").Value = 924.7999
$ws.Range("N138").Value = -16540.4645
$ws.Range("H141").Value = 7519.0557
$ws.Range("I141").Value = 8987.691999999999
$ws.Range("J141").Value = 3700.6
$ws.Range("K141").Value = 26963.076
$ws.Range("L141").Value = 11101.8
$ws.Range("M141").Value = -21783.076
$ws.Range("N141").Value = -21461.8

# --- ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H74").Value = 855.55554
$ws.Range("I74").Value = 846.2353000000001
$ws.Range("K74").Value = 846.2353000000001
$ws.Range("M74").Value = 27.76469999999995
$ws.Range("H77").Value = 855.55554
$ws.Range("I77").Value = 846.2353000000001
$ws.Range("K77").Value = 4231.1765
$ws.Range("M77").Value = 136.8234999999995
$ws.Range("H102").Value = 27795428
$ws.Range("I102").Value = 33354226
$ws.Range("J102").Value = 1445
$ws.Range("K102").Value = 33354226
$ws.Range("L102").Value = 1445
$ws.Range("M102").Value = -33352604
$ws.Range("N102").Value = -4689
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H132").Value = 2765.16
$ws.Range("I132").Value = 3075.5334
$ws.Range("J132").Value = 2299.6
$ws.Range("K132").Value = 9226.600199999999
$ws.Range("L132").Value = 6898.799999999999
$ws.Range("M132").Value = -6696.600199999999
$ws.Range("N132").Value = -11958.8
$ws.Range("H139").Value = 50607
$ws.Range("J139").Value = 50607
$ws.Range("L139").Value = 50607
$ws.Range("N139").Value = -60887

# --- BSM ---
$ws = $wb.Worksheets("BSM")
$ws.Range("H105").Value = 201980320
$ws.Range("I105").Value = 201980320
$ws.Range("K105").Value = 201980320
$ws.Range("M105").Value = -201978573

# --- CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H99").Value = 1909.25
$ws.Range("I99").Value = 1792
$ws.Range("J99").Value = 2104.6667
$ws.Range("K99").Value = 1792
$ws.Range("L99").Value = 2104.6667
$ws.Range("M99").Value = -294
$ws.Range("N99").Value = -5100.6667
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H126").Value = 1909.25
$ws.Range("I126").Value = 1792
$ws.Range("J126").Value = 2104.6667
$ws.Range("K126").Value = 5376
$ws.Range("L126").Value = 6314.000100000001
$ws.Range("M126").Value = -2906
$ws.Range("N126").Value = -11254.0001
$ws.Range("H134").Value = 11495772
$ws.Range("I134").Value = 14494109
$ws.Range("K134").Value = 43482327
$ws.Range("M134").Value = -43479792

# --- CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H16").Value = 2901
$ws.Range("J16").Value = 2901
$ws.Range("L16").Value = 8703
$ws.Range("N16").Value = -9049
$ws.Range("H68").Value = 1368.2778
$ws.Range("I68").Value = 748.2
$ws.Range("J68").Value = 1606.7693
$ws.Range("K68").Value = 2244.6
$ws.Range("L68").Value = 4820.3079
$ws.Range("M68").Value = -1433.6
$ws.Range("N68").Value = -6442.3079
$ws.Range("H71").Value = 1368.2778
$ws.Range("I71").Value = 748.2
$ws.Range("J71").Value = 1606.7693
$ws.Range("K71").Value = 6733.8
$ws.Range("L71").Value = 14460.9237
$ws.Range("M71").Value = -2677.8
$ws.Range("N71").Value = -22572.9237
$ws.Range("H94").Value = 5240
$ws.Range("J94").Value = 5322.222
$ws.Range("L94").Value = 15966.666
$ws.Range("N94").Value = -17318.666
$ws.Range("H98").Value = 725.5
$ws.Range("I98").Value = 146
$ws.Range("K98").Value = 438
$ws.Range("M98").Value = 1060
$ws.Range("H131").Value = 27028602
$ws.Range("I131").Value = 200000860
$ws.Range("J131").Value = 1684.9375
$ws.Range("K131").Value = 600002580
$ws.Range("L131").Value = 5054.8125
$ws.Range("M131").Value = -599997540
$ws.Range("N131").Value = -15134.8125
$ws.Range("H133").Value = 3379.88
$ws.Range("J133").Value = 3568.087
$ws.Range("L133").Value = 10704.261
$ws.Range("N133").Value = -20824.261
$ws.Range("H134").Value = 3740.4736
$ws.Range("I134").Value = 2511.5
$ws.Range("J134").Value = 4307.6924
$ws.Range("K134").Value = 7534.5
$ws.Range("L134").Value = 12923.0772
$ws.Range("M134").Value = -2464.5
$ws.Range("N134").Value = -23063.0772
$ws.Range("H137").Value = 5403.609
$ws.Range("I137").Value = 816.6667
$ws.Range("J137").Value = 6091.65
$ws.Range("K137").Value = 2450.0001
$ws.Range("L137").Value = 18274.95
$ws.Range("M137").Value = 2649.9999
$ws.Range("N137").Value = -28474.95

# --- GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H70").Value = 75005200
$ws.Range("I70").Value = 250000000
$ws.Range("J70").Value = 40006240
$ws.Range("K70").Value = 250000000
$ws.Range("L70").Value = 40006240
$ws.Range("M70").Value = -249999730
$ws.Range("N70").Value = -40006780
$ws.Range("H73").Value = 75005200
$ws.Range("I73").Value = 250000000
$ws.Range("J73").Value = 40006240
$ws.Range("K73").Value = 250000000
$ws.Range("L73").Value = 40006240
$ws.Range("M73").Value = -249999064
$ws.Range("N73").Value = -40008112
$ws.Range("H102").Value = 1433.3871
$ws.Range("I102").Value = 1405.0416
$ws.Range("J102").Value = 1530.5714
$ws.Range("K102").Value = 1405.0416
$ws.Range("L102").Value = 1530.5714
$ws.Range("M102").Value = 216.9584
$ws.Range("N102").Value = -4774.5714

# --- LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 2222
$ws.Range("I7").Value = 2222
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2222
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2110
$ws.Range("N7").ClearContents()
$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 900
$ws.Range("K12").Value = 900
$ws.Range("M12").Value = -730
$ws.Range("H22").Value = 913.0909
$ws.Range("I22").Value = 409
$ws.Range("K22").Value = 409
$ws.Range("M22").Value = -114
$ws.Range("H27").Value = 913.0909
$ws.Range("I27").Value = 409
$ws.Range("K27").Value = 409
$ws.Range("M27").Value = -302
$ws.Range("H53").Value = 4833.3335
$ws.Range("J53").Value = 7000
$ws.Range("L53").Value = 7000
$ws.Range("N53").Value = -8036
$ws.Range("H126").Value = 2222
$ws.Range("I126").Value = 2222
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6666
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4196
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 65150.438
$ws.Range("I132").Value = 2299.889
$ws.Range("J132").Value = 145958.28
$ws.Range("K132").Value = 6899.667
$ws.Range("L132").Value = 437874.84
$ws.Range("M132").Value = -4369.667
$ws.Range("N132").Value = -442934.84
$ws.Range("H136").Value = 2481.6
$ws.Range("I136").Value = 2601.375
$ws.Range("J136").Value = 2002.5
$ws.Range("K136").Value = 7804.125
$ws.Range("L136").Value = 6007.5
$ws.Range("M136").Value = -5254.125
$ws.Range("N136").Value = -11107.5
